$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.007741359813101
$ws.Range("D2").Value = 1.010450795242707
$ws.Range("E2").Value = 1.010378232713478
$ws.Range("F2").Value = 1.005934661921953
$ws.Range("J2").Value = 1.013011236820697
$ws.Range("K2").Value = 1.013320948185163
$ws.Range("L2").Value = 1.013248605009605
$ws.Range("M2").Value = 1.008818532729725
$ws.Range("N2").Value = 1.008424484006033
$ws.Range("C3").Value = 1.008885489017892
$ws.Range("D3").Value = 1.011447612279555
$ws.Range("E3").Value = 1.011354631784088
$ws.Range("F3").Value = 1.007732031210571
$ws.Range("J3").Value = 1.013785937500267
$ws.Range("K3").Value = 1.014121302574718
$ws.Range("L3").Value = 1.014028582638701
$ws.Range("M3").Value = 1.010416174298673
$ws.Range("N3").Value = 1.008689512514675
$ws.Range("C4").Value = 1.00962571700567
$ws.Range("D4").Value = 1.012092817872978
$ws.Range("E4").Value = 1.01198666391706
$ws.Range("F4").Value = 1.008894589360257
$ws.Range("J4").Value = 1.014286659102448
$ws.Range("K4").Value = 1.014638783828894
$ws.Range("L4").Value = 1.014532912232917
$ws.Range("M4").Value = 1.011449090507858
$ws.Range("N4").Value = 1.008860573501175
$ws.Range("C5").Value = 1.009936887534993
$ws.Range("D5").Value = 1.012364111420605
$ws.Range("E5").Value = 1.012252428344765
$ws.Range("F5").Value = 1.009383228894149
$ws.Range("J5").Value = 1.014497029797006
$ws.Range("K5").Value = 1.014856238440871
$ws.Range("L5").Value = 1.014744845756573
$ws.Range("M5").Value = 1.011883130412383
$ws.Range("N5").Value = 1.008932384960895
$ws.Range("C6").Value = 1.009989133262976
$ws.Range("D6").Value = 1.012409665728677
$ws.Range("E6").Value = 1.012297054812175
$ws.Range("F6").Value = 1.009465268058941
$ws.Range("J6").Value = 1.014532344233848
$ws.Range("K6").Value = 1.014892744532218
$ws.Range("L6").Value = 1.014780425290232
$ws.Range("M6").Value = 1.01195599626664
$ws.Range("N6").Value = 1.008944436417451
$ws.Range("C7").Value = 1.009629874961979
$ws.Range("D7").Value = 1.012096442715623
$ws.Range("E7").Value = 1.011990214844175
$ws.Range("F7").Value = 1.008901118965581
$ws.Range("J7").Value = 1.014289470605075
$ws.Range("K7").Value = 1.014641689837406
$ws.Range("L7").Value = 1.014535744437633
$ws.Range("M7").Value = 1.011454890937521
$ws.Range("N7").Value = 1.008861533451982
$ws.Range("C8").Value = 1.008128045961337
$ws.Range("D8").Value = 1.010787634154177
$ws.Range("E8").Value = 1.010708163237112
$ws.Range("F8").Value = 1.006542192680028
$ws.Range("J8").Value = 1.013273167878997
$ws.Range("K8").Value = 1.013591515823635
$ws.Range("L8").Value = 1.013512279211207
$ws.Range("M8").Value = 1.009358647616566
$ws.Range("N8").Value = 1.008514140949303
$ws.Range("C9").Value = 1.005480745576251
$ws.Range("D9").Value = 1.008482791098277
$ws.Range("E9").Value = 1.00845077064082
$ws.Range("F9").Value = 1.00238148300482
$ws.Range("J9").Value = 1.01147793231855
$ws.Range("K9").Value = 1.011737826317246
$ws.Range("L9").Value = 1.011705916681196
$ws.Range("M9").Value = 1.005657775411462
$ws.Range("N9").Value = 1.007898677791794
$ws.Range("C10").Value = 1.003715114672957
$ws.Range("D10").Value = 1.006947094622776
$ws.Range("E10").Value = 1.006946921126704
$ws.Range("F10").Value = 0.9996043609359123
$ws.Range("J10").Value = 1.010278050462231
$ws.Range("K10").Value = 1.010499811863726
$ws.Range("L10").Value = 1.010499639027785
$ws.Range("M10").Value = 1.003185244096015
$ws.Range("N10").Value = 1.00748611456252
$ws.Range("C11").Value = 1.002950357702456
$ws.Range("D11").Value = 1.006282302090868
$ws.Range("E11").Value = 1.006295971652534
$ws.Range("F11").Value = 0.998400907678803
$ws.Range("J11").Value = 1.009757737288748
$ws.Range("K11").Value = 1.009963187190324
$ws.Range("L11").Value = 1.009976802697978
$ws.Range("M11").Value = 1.002113232761298
$ws.Range("N11").Value = 1.007306929117971
$ws.Range("C12").Value = 1.002666255427317
$ws.Range("D12").Value = 1.006035392536111
$ws.Range("E12").Value = 1.006054212209226
$ws.Range("F12").Value = 0.9979537378880738
$ws.Range("J12").Value = 1.009564353877718
$ws.Range("K12").Value = 1.0097637751944
$ws.Range("L12").Value = 1.009782519401926
$ws.Range("M12").Value = 1.001714820554971
$ws.Range("N12").Value = 1.007240289530172
$ws.Range("C13").Value = 1.002727198080188
$ws.Range("D13").Value = 1.006088354423225
$ws.Range("E13").Value = 1.006106069006355
$ws.Range("F13").Value = 0.9980496644779331
$ws.Range("J13").Value = 1.009605840558904
$ws.Range("K13").Value = 1.009806553656861
$ws.Range("L13").Value = 1.009824197417062
$ws.Range("M13").Value = 1.0018002914309
$ws.Range("N13").Value = 1.007254587674138
$ws.Range("C14").Value = 1.002926874503961
$ws.Range("D14").Value = 1.006261891993281
$ws.Range("E14").Value = 1.006275987096645
$ws.Range("F14").Value = 0.9983639476905005
$ws.Range("J14").Value = 1.009741754533931
$ws.Range("K14").Value = 1.009946705493142
$ws.Range("L14").Value = 1.009960744795007
$ws.Range("M14").Value = 1.002080304423726
$ws.Range("N14").Value = 1.007301422347768
$ws.Range("C15").Value = 1.003049896670782
$ws.Range("D15").Value = 1.006368817331358
$ws.Range("E15").Value = 1.006380683444272
$ws.Range("F15").Value = 0.9985575672192483
$ws.Range("J15").Value = 1.009825480188881
$ws.Range("K15").Value = 1.010033046245151
$ws.Range("L15").Value = 1.01004486566652
$ws.Range("M15").Value = 1.002252800264958
$ws.Range("N15").Value = 1.007330267834748
$ws.Range("C16").Value = 1.003765864002442
$ws.Range("D16").Value = 1.00699121819798
$ws.Range("E16").Value = 1.00699012712017
$ws.Range("F16").Value = 0.9996842092912568
$ws.Range("J16").Value = 1.010312565802186
$ws.Range("K16").Value = 1.010535413975111
$ws.Range("L16").Value = 1.01053432703791
$ws.Range("M16").Value = 1.003256359774522
$ws.Range("N16").Value = 1.007497995023939
$ws.Range("C17").Value = 1.004214908215161
$ws.Range("D17").Value = 1.007381678918184
$ws.Range("E17").Value = 1.007372474500552
$ws.Range("F17").Value = 1.00039066149902
$ws.Range("J17").Value = 1.01061789761435
$ws.Range("K17").Value = 1.01085038555817
$ws.Range("L17").Value = 1.010841215264662
$ws.Range("M17").Value = 1.003885486680965
$ws.Range("N17").Value = 1.007603060116528
$ws.Range("C18").Value = 1.004476806420504
$ws.Range("D18").Value = 1.007609445004532
$ws.Range("E18").Value = 1.007595513323928
$ws.Range("F18").Value = 1.000802633153154
$ws.Range("J18").Value = 1.010795919667981
$ws.Range("K18").Value = 1.011034049629468
$ws.Range("L18").Value = 1.011020168887485
$ws.Range("M18").Value = 1.004252312479569
$ws.Range("N18").Value = 1.007664290470718
$ws.Range("C19").Value = 1.004566103363045
$ws.Range("D19").Value = 1.007687110294921
$ws.Range("E19").Value = 1.007671567600543
$ws.Range("F19").Value = 1.000943089951269
$ws.Range("J19").Value = 1.01085660830383
$ws.Range("K19").Value = 1.011096665316714
$ws.Range("L19").Value = 1.011081179180929
$ws.Range("M19").Value = 1.00437736831558
$ws.Range("N19").Value = 1.007685159579239
$ws.Range("C20").Value = 1.004166732289501
$ws.Range("D20").Value = 1.007339784441221
$ws.Range("E20").Value = 1.007331449978018
$ws.Range("F20").Value = 1.000314875259818
$ws.Range("J20").Value = 1.010585145950746
$ws.Range("K20").Value = 1.01081659764068
$ws.Range("L20").Value = 1.010808294154659
$ws.Range("M20").Value = 1.003818001184112
$ws.Range("N20").Value = 1.007591793044343
$ws.Range("C21").Value = 1.002868075835822
$ws.Range("D21").Value = 1.00621078885949
$ws.Range("E21").Value = 1.00622594958189
$ws.Range("F21").Value = 0.9982714034586953
$ws.Range("J21").Value = 1.009701734474589
$ws.Range("K21").Value = 1.009905436653439
$ws.Range("L21").Value = 1.009920537181223
$ws.Range("M21").Value = 1.001997853721787
$ws.Range("N21").Value = 1.007287632975731
$ws.Range("C22").Value = 1.002051340423527
$ws.Range("D22").Value = 1.005501082244604
$ws.Range("E22").Value = 1.00553106277097
$ws.Range("F22").Value = 0.9969856964212362
$ws.Range("J22").Value = 1.009145626728355
$ws.Range("K22").Value = 1.0093320566831
$ws.Range("L22").Value = 1.009361913377761
$ws.Range("M22").Value = 1.000852179619299
$ws.Range("N22").Value = 1.007095919966662
$ws.Range("C23").Value = 1.002484329161242
$ws.Range("D23").Value = 1.005877298825826
$ws.Range("E23").Value = 1.00589941852797
$ws.Range("F23").Value = 0.9976673630241024
$ws.Range("J23").Value = 1.009440494421733
$ws.Range("K23").Value = 1.009636064068379
$ws.Range("L23").Value = 1.009658094284072
$ws.Range("M23").Value = 1.001459647627041
$ws.Range("N23").Value = 1.007197595927036
$ws.Range("C24").Value = 1.004188500983482
$ws.Range("D24").Value = 1.007358714699621
$ws.Range("E24").Value = 1.007349987125685
$ws.Range("F24").Value = 1.000349120075985
$ws.Range("J24").Value = 1.010599945242686
$ws.Range("K24").Value = 1.010831865111529
$ws.Range("L24").Value = 1.010823169938823
$ws.Range("M24").Value = 1.003848495386097
$ws.Range("N24").Value = 1.00759688431099
$ws.Range("C25").Value = 1.006165258446153
$ws.Range("D25").Value = 1.009078487903275
$ws.Range("E25").Value = 1.009034163554825
$ws.Range("F25").Value = 1.003457663494369
$ws.Range("J25").Value = 1.011942574428573
$ws.Range("K25").Value = 1.01221743333349
$ws.Range("L25").Value = 1.012173256498776
$ws.Range("M25").Value = 1.006615429486483
$ws.Range("N25").Value = 1.008058185146165
